# Update column F (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -2
    3  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = -4
    8  = -1
    9  = -1
    10 = -1
    11 = -4
    12 = 3
    13 = 1
    14 = 6
    15 = 6
    16 = -6
    18 = -2
    20 = -1
    21 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
